$p = $ppt.ActivePresentation

# Locate the slide/shape containing the "according to ISO 9126." sentence.
$targetSlide = $null
$targetShape = $null
foreach ($s in $p.Slides) {
    foreach ($shp in $s.Shapes) {
        if ($shp.HasTextFrame) {
            $t = $shp.TextFrame.TextRange.Text
            if ($t.Contains("according to ISO 9126.")) {
                $targetSlide = $s
                $targetShape = $shp
            }
        }
    }
}

$tr = $targetShape.TextFrame.TextRange

# Step 1: the existing run "according to ISO 9126." gains a trailing space,
# without touching its formatting (same run, same rPr).
$full = $tr.Text
$searchText = "according to ISO 9126."
$idx = $full.IndexOf($searchText)
$startPos = $idx + 1
$sub = $tr.Characters($startPos, $searchText.Length)
$sub.Text = "according to ISO 9126. "

# Step 2: append a brand-new run "(ACSTC)" right after it, inheriting the
# bold/black "Inter Bold" character formatting of the text it follows.
$tr2 = $targetShape.TextFrame.TextRange
$tr2.InsertAfter("(ACSTC)") | Out-Null
